# Finalized code and info box before exe generation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report cutoff date text (shared string, in merged cell A2:F2)
$ws.Range("A2").Value = "REPORT - CUTOFF DATE: 2024-10-01"

# Materialize the previously-empty info box row (A1:F1, already merged &
# included in the sheet dimension) by copying the formatting down from the
# row below it - this stamps cells A1:F1 with the same style (s="1") without
# giving them any value/content.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A1:F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the CURRENT EXPENDITURES (column C) figures that were still at
# their placeholder 0 value.
$ws.Range("C4").Value = 4706.24
$ws.Range("C5").Value = 14.17
$ws.Range("C6").Value = 419.88
$ws.Range("C8").Value = 5569.36
$ws.Range("C10").Value = -675.2
$ws.Range("C11").Value = 460.48
